$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "StatQuery" Cypher text (replaces the old aggregate-counts query used
# by the CasesTab / SamplesTab / FilesTab rows in column C).
$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Boston Terrier']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value2 = $newStatQuery
$ws.Range("C3").Value2 = $newStatQuery
$ws.Range("C4").Value2 = $newStatQuery

# Window/view tidy-up: drop the old 160% zoom + B1 top-left scroll position,
# reset to 100%, and move the selection to C4 (matches the saved view state).
$excel.ActiveWindow.Zoom = 100
$ws.Range("C4").Select() | Out-Null
